# Adding a function to compare h to optimum
# Renumber the grouping id in column A of the "OutAssignment" sheet for
# rows 17-63: the blocks that were labeled 13,5,6,7,8 become 5,6,7,8,9
# respectively (i.e. shifted down so they continue sequentially from the
# earlier rows in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OutAssignment")

$ws.Range("A17:A37").Value = 5
$ws.Range("A38:A46").Value = 6
$ws.Range("A47:A51").Value = 7
$ws.Range("A52:A60").Value = 8
$ws.Range("A61:A63").Value = 9
